$d = $word.ActiveDocument

# The document's first paragraph currently reads:  <page>126v</page>
# (as three runs: a gray "<page>" tag, a black "126v", and a gray
# "</page>" tag). We need to insert a brand new run containing the
# text "tirm" immediately in front of the "<page>" run (i.e. right at
# the very start of the paragraph / document), styled the same way as
# the other gray tag markers: Courier New, 9pt (sz/szCs=18 half-pts),
# color a9a9a9.

$target = $d.Paragraphs(1).Range.Duplicate
$target.Collapse(1)   # wdCollapseStart -> put the insertion point before "<page>"

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body><w:p><w:r>' +
       '<w:rPr>' +
       '<w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/>' +
       '<w:color w:val="a9a9a9"/>' +
       '<w:sz w:val="18"/>' +
       '<w:szCs w:val="18"/>' +
       '<w:rtl w:val="0"/>' +
       '</w:rPr>' +
       '<w:t xml:space="preserve">tirm</w:t>' +
       '</w:r></w:p></w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
